$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.361.48"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "1.895.37"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.55"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -3.19%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4775"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +2.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4058"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -1.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08049"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("E10").Value = "  -1.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.29"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +4.04%  "
$ws.Range("D12").Value = "1.879.08"
$ws.Range("E12").Value = "  -1.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.940"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.061"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.66"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06666"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +1.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001029"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").Value = "29.384.88"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.522"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.68"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.156"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -2.01%  "
$ws.Range("D25").Value = "2.157.54"
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.40"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.74"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.033"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +5.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.086"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -2.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.05"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.021"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -5.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09491"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.528"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.385"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -3.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.369"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02246"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -1.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06046"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.171"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -0.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5856"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.864"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -6.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1840"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.12"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.419"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +3.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.288"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +3.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07715"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +2.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.16"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5497"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -1.59%  "
$ws.Range("E48").Value = "  -0.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.83"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.2953"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.64"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -1.12%  "
